$d = $word.ActiveDocument

# --- 1. Remove the stray "_GoBack" bookmark that sits after "Oracle D" ---
# (Word auto-tracks the last edit location with a hidden "_GoBack" bookmark;
#  it gets re-created near the new edit below, but make sure the old one
#  is gone explicitly too.)
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 2. Split "...Teradata, SQL to bring..." into three runs, inserting "and " ---
$rng = $d.Content
$found = $rng.Find.Execute("Teradata, ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Collapse to the point right before "SQL" and type the new word.
    $rng.Collapse(0)
    $rng.InsertAfter("and ")

    # Toggling a character property and reverting it forces the newly
    # inserted "and " text to live in its own run instead of being
    # silently re-merged into the identically-formatted neighbour run.
    $rng.Bold = 1
    $rng.Bold = 0

    # --- 3. Re-insert the "_GoBack" bookmark right after "and ", before "SQL" ---
    $rng.Collapse(0)
    $d.Bookmarks.Add("_GoBack", $rng)
}
